$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" column header in F1
$ws.Range("F1").Value = "time_taken"

# Populate F2:F6 with the time_taken values (plain text, like the other data columns)
$ws.Range("F2").Value = "2021-10-05 13:40:48.930303"
$ws.Range("F3").Value = "2021-10-05 13:40:48.930315"
$ws.Range("F4").Value = "2021-10-05 13:40:48.930318"
$ws.Range("F5").Value = "2021-10-05 13:40:48.930321"
$ws.Range("F6").Value = "2021-10-05 13:40:48.930324"

# Match the header formatting (bold, bordered, centered) used by the other header cells
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
